$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Drop the stale "_GoBack" bookmark that used to sit right after the
#    "2. Исходные данные к проекту:" heading run.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2) Re-colour the five "обеспечивать ..." bullet items (driver ads,
#    passenger ads, authorization, registration, view-list) from green
#    to cyan highlighting. wdColorIndex 3 == wdTurquoise (OOXML "cyan").
# ---------------------------------------------------------------------
$wdTurquoise = 3

$paraTexts = @(
    "обеспечивать создание объявлений в качестве водителя;",
    "обеспечивать создание объявлений в качестве попутчика;",
    "обеспечивать авторизацию пользователей;",
    "обеспечивать регистрацию пользователей;",
    "обеспечивать просмотр списка объявлений;"
)

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text
    foreach ($needle in $paraTexts) {
        if ($text.StartsWith($needle)) {
            $p.Range.HighlightColorIndex = $wdTurquoise
        }
    }
}

# ---------------------------------------------------------------------
# 3) Re-insert "_GoBack" inside the (now cyan) driver-ads bullet, right
#    in the middle of "водителя" -> "вод" | "ителя", matching the split
#    the author made while editing that line.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("в качестве вод", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertPos = $rng.End
$insertRng = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $insertRng) | Out-Null
